$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continuation of the participant table: extend the placeholder rows from
# 10 participants to 23 participants (rows 11..24), following the exact
# same pattern as the existing rows (3..10): column A holds the
# participant index (row - 1), columns B..P are blank placeholders
# waiting to be filled in.
$columns = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")

for ($row = 11; $row -le 24; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 1

    foreach ($col in $columns) {
        $ws.Range("$col$row").Value = ""
    }
}
